$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The sheet described a set of "Transportation Personal / Transportation
# Freight" blend rows; the new data re-organises this into "Fuel Blends"
# (Gasoline_Transportation / Diesel_Transportation) rows, adds a class-min
# row for each fuel, a Renewable Gasoline / Renewable Diesel technology row,
# and shifts/updates several formulas. Easiest/most-robust approach: clear
# the previous data block and rewrite the new block from scratch.
# ---------------------------------------------------------------------------

$ws.Range("A3:X8").ClearContents() | Out-Null

# --- Row 3: Gasoline_Transportation / Market share_class_min --------------
$ws.Range("A3").Value = "CIMS.CAN.ON.Fuel Blends.Gasoline_Transportation"
$ws.Range("B3").Value = "Service"
$ws.Range("C3").Value = "ON"
$ws.Range("E3").Value = "Gasoline_Transportation"
$ws.Range("G3").Value = "Market share_class_min"
$ws.Range("H3").Value = "Gasoline_Transportation"

$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Formula = "=0.67*0.05/(0.67*0.05+1*0.95)"
$ws.Range("Q3").Formula = "=0.67*0.1/(0.67*0.1+1*0.9)"
$ws.Range("R3").Formula = "=0.67*0.11/(0.67*0.11+1*0.89)"
$ws.Range("S3").Formula = "=0.67*0.15/(0.67*0.15+1*0.85)"
$ws.Range("T3:W3").Formula = "=S3"

$ws.Range("X3").Value = "2011-2015 (?) federl standard; 2016 (?) to 2030 Ontario standard"

# --- Row 4: Gasoline_Transportation / Ethanol / Market share_class --------
$ws.Range("A4").Value = "CIMS.CAN.ON.Fuel Blends.Gasoline_Transportation"
$ws.Range("B4").Value = "Service"
$ws.Range("C4").Value = "ON"
$ws.Range("E4").Value = "Gasoline_Transportation"
$ws.Range("F4").Value = "Ethanol"
$ws.Range("G4").Value = "Market share_class"
$ws.Range("H4").Value = "Gasoline_Transportation"

# --- Row 5: Gasoline_Transportation / Renewable Gasoline / Market share_class
$ws.Range("A5").Value = "CIMS.CAN.ON.Fuel Blends.Gasoline_Transportation"
$ws.Range("B5").Value = "Service"
$ws.Range("C5").Value = "ON"
$ws.Range("E5").Value = "Gasoline_Transportation"
$ws.Range("F5").Value = "Renewable Gasoline"
$ws.Range("G5").Value = "Market share_class"
$ws.Range("H5").Value = "Gasoline_Transportation"

# --- Row 6: Diesel_Transportation / Market share new_min -------------------
$ws.Range("A6").Value = "CIMS.CAN.ON.Fuel Blends.Diesel_Transportation"
$ws.Range("B6").Value = "Service"
$ws.Range("C6").Value = "ON"
$ws.Range("E6").Value = "Diesel_Transportation"
$ws.Range("G6").Value = "Market share new_min"
$ws.Range("H6").Value = "Diesel_Transportation"

$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0
$ws.Range("O6").Value = 0
$ws.Range("P6").Formula = "=0.87*0.02/(0.87*0.02+1*0.98)"
$ws.Range("Q6").Formula = "=0.87*0.04/(0.87*0.04+1*0.96)"
$ws.Range("R6:W6").Formula = "=Q6"

$ws.Range("X6").Value = "2011-2015 (?) federl standard; 2016 (?) to 2030 Ontario standard"

# --- Row 7: Diesel_Transportation / Biodiesel / Market share_class --------
$ws.Range("A7").Value = "CIMS.CAN.ON.Fuel Blends.Diesel_Transportation"
$ws.Range("B7").Value = "Service"
$ws.Range("C7").Value = "ON"
$ws.Range("E7").Value = "Diesel_Transportation"
$ws.Range("F7").Value = "Biodiesel"
$ws.Range("G7").Value = "Market share_class"
$ws.Range("H7").Value = "Diesel_Transportation"

# --- Row 8: Diesel_Transportation / Renewable Diesel / Market share_class -
$ws.Range("A8").Value = "CIMS.CAN.ON.Fuel Blends.Diesel_Transportation"
$ws.Range("B8").Value = "Service"
$ws.Range("C8").Value = "ON"
$ws.Range("E8").Value = "Diesel_Transportation"
$ws.Range("F8").Value = "Renewable Diesel"
$ws.Range("G8").Value = "Market share_class"
$ws.Range("H8").Value = "Diesel_Transportation"

$ws.Range("A1:X8").Select() | Out-Null
